# Round the numeric estimate/std.error/statistic/p.value columns (B:E, rows 2-17)
# of the model-results table down to 2 decimal places (manipulation check /
# promotion-based final results / exploration post-hoc update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09
$ws.Range("C2").Value = 0.77
$ws.Range("D2").Value = 0.11
$ws.Range("E2").Value = 0.91

$ws.Range("B3").Value = -0.01
$ws.Range("C3").Value = 0.01
$ws.Range("D3").Value = -0.49
$ws.Range("E3").Value = 0.63

$ws.Range("B4").Value = -0.2
$ws.Range("C4").Value = 0.2
$ws.Range("D4").Value = -0.97
$ws.Range("E4").Value = 0.33

$ws.Range("B5").Value = 0.05
$ws.Range("C5").Value = 0.06
$ws.Range("D5").Value = 0.75
$ws.Range("E5").Value = 0.45

$ws.Range("B6").Value = -0.06
$ws.Range("C6").Value = 0.16
$ws.Range("D6").Value = -0.39

$ws.Range("B7").Value = -0.01
$ws.Range("C7").Value = 0.09
$ws.Range("D7").Value = -0.08
$ws.Range("E7").Value = 0.93

$ws.Range("B8").Value = 0.04
$ws.Range("D8").Value = 0.47
$ws.Range("E8").Value = 0.64

$ws.Range("B9").Value = 0.29
$ws.Range("C9").Value = 0.2
$ws.Range("D9").Value = 1.46
$ws.Range("E9").Value = 0.14

$ws.Range("B10").Value = -0.46
$ws.Range("C10").Value = 0.2
$ws.Range("D10").Value = -2.24
$ws.Range("E10").Value = 0.02

$ws.Range("B11").Value = -0.03
$ws.Range("C11").Value = 0.23
$ws.Range("E11").Value = 0.9

$ws.Range("B12").Value = 0.31
$ws.Range("C12").Value = 0.28
$ws.Range("D12").Value = 1.1
$ws.Range("E12").Value = 0.27

$ws.Range("B13").Value = 0.0
$ws.Range("C13").Value = 0.28
$ws.Range("D13").Value = 0.01
$ws.Range("E13").Value = 0.99

$ws.Range("B14").Value = 0.38
$ws.Range("C14").Value = 0.29

$ws.Range("B15").Value = 0.23
$ws.Range("C15").Value = 0.3
$ws.Range("D15").Value = 0.74
$ws.Range("E15").Value = 0.46

$ws.Range("B16").Value = -0.12
$ws.Range("C16").Value = 0.28
$ws.Range("E16").Value = 0.65

$ws.Range("B17").Value = 0.22
$ws.Range("C17").Value = 0.3
